# "création vue initialisation projet"
# Update student "Numero" (col A) to reflect promotion 2017 (was 2015)
# and refresh "Moyenne de l'etudiant" (col E) values for the updated cohort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: every student number for rows 3..63 moves from the 2015xxxx
# series to the matching 2017xxxx series (same last 4 digits, +20000).
for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + 20000
    }
}

# Column E: new "moyenne" values for the rows where it changed.
$newAverages = @{
    3 = 19; 4 = 7; 5 = 15; 6 = 17; 7 = 8; 8 = 16; 9 = 19; 10 = 9; 11 = 10; 12 = 19;
    14 = 15; 15 = 5; 16 = 14; 17 = 9; 18 = 9; 19 = 13; 20 = 9; 21 = 17; 22 = 15; 23 = 17; 24 = 9; 25 = 18;
    27 = 18; 28 = 8; 29 = 12; 30 = 8; 31 = 6; 32 = 20; 33 = 10; 34 = 18; 35 = 6;
    37 = 20; 38 = 7; 39 = 5; 40 = 5; 41 = 8; 42 = 12; 43 = 13; 44 = 13; 45 = 14; 46 = 6; 47 = 18; 48 = 19; 49 = 16; 50 = 18; 51 = 19; 52 = 19; 53 = 17; 54 = 12; 55 = 8;
    57 = 17; 58 = 5; 59 = 11; 60 = 11; 61 = 18; 62 = 9; 63 = 10
}

foreach ($row in $newAverages.Keys) {
    $ws.Cells.Item($row, 5).Value = $newAverages[$row]
}
